$d = $word.ActiveDocument

# Locate the stretch of text that needs to be restructured: the old
# "Random-Forest-Regressionsmodells als Grundlage ... Streichmaschinenanlage."
# tail of the sentence.
$target = $d.Content
$found = $target.Find.Execute("Random-Forest-Regressionsmodells als Grundlage einer teilautonomen und zukünftig vollautonomen sowie gemäß EU AI Act transparenten und sicheren Prozesssteuerung an einer Streichmaschinenanlage.")

if ($found) {
    # Replace it with the same text split into multiple runs: the new
    # English term "Machine" is wrapped in proofErr spellStart/spellEnd
    # (as Word marks an unrecognized word), followed by "-Learning-Modells",
    # a standalone space run, and finally the remainder of the sentence.
    $ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
    $xml = '<w:p ' + $ns + '>' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Machine</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>-Learning-Modells</w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>als Grundlage einer teilautonomen und zukünftig vollautonomen sowie gemäß EU AI Act transparenten und sicheren Prozesssteuerung an einer Streichmaschinenanlage.</w:t></w:r>' +
           '</w:p>'
    $target.InsertXML($xml)
}

Write-Output "done"
